$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (N1, O1, P1), copying the style from the
# existing header cell M1 so the new headers get the same bold/border
# formatting (style index) as the rest of row 1.
$ws.Range("M1").Copy()
$ws.Range("N1:P1").PasteSpecial(-4122)
$ws.Range("N1").Value = "renewd"
$ws.Range("O1").Value = "PlanID"
$ws.Range("P1").Value = "iteration"

# Update every data row (2-103): rename the status tag in column D from
# "215_2" (with leading BOM) to "215_2n", and populate the three new
# columns N (renewd), O (PlanID), P (iteration).
for ($r = 2; $r -le 103; $r++) {
    $ws.Range("D$r").Value = "﻿215_2n"
    $ws.Range("N$r").Value = "after"
    $ws.Range("O$r").Value = "502-0147678"
    $ws.Range("P$r").Value = 14
}
